# feat: add 2022-Q1 data
#
# - Renames the old "总计" (summary) sheet into a new "2022-Q1" fund-holdings
#   sheet (same per-quarter layout used by 2020-Q4..2021-Q4).
# - Appends a brand new "总计" sheet at the end, containing the original
#   summary rows plus a new leading "2022-Q1" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: the existing "总计" sheet (index 6) becomes "2022-Q1" and is
# repopulated with the quarterly fund-holdings table.
# ---------------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item(6)
$q1Sheet.Name = "2022-Q1"

# The previous "总计" sheet only used A1:D6; the new table is a strict
# superset (A1:H13), so overwriting header/data cells directly below is
# enough - nothing needs clearing first, and clearing would also wipe the
# cell styling we still need to read for the header/index cells further
# down.

$q1Headers = @('基金代码','基金名称','基金规模','股票总仓位','仓位占比','持有市值(亿元)','仓位排名')
for ($c = 0; $c -lt $q1Headers.Length; $c++) {
    $cell = $q1Sheet.Cells.Item(1, 2 + $c)
    $cell.Value = $q1Headers[$c]
}

$rows2022Q1 = @(
    @(0,'420005','天弘周期策略混合','5.25','89.31','5.25','0.2756',7),
    @(1,'420001','天弘精选混合','7.16','71.80','3.25','0.2327',5),
    @(2,'007202','天弘优质成长企业精选混合','4.81','92.52','4.27','0.2054',9),
    @(3,'001030','天弘云端生活优选灵活配置混合','1.61','79.35','5.99','0.0964',4),
    @(4,'004694','天弘策略精选灵活配置混合A','1.11','80.93','3.89','0.0432',4),
    @(5,'009186','天弘聚新三个月定期开放混合A','2.93','23.24','0.99','0.0290',9),
    @(6,'002388','天弘裕利灵活配置混合A','1.87','29.63','1.53','0.0286',9),
    @(7,'010331','天弘消费股票A','0.47','83.48','5.49','0.0258',4),
    @(8,'005997','天弘裕利灵活配置混合C','0.99','29.63','1.53','0.0151',9),
    @(9,'010332','天弘消费股票C','0.20','83.48','5.49','0.0110',4),
    @(10,'004748','天弘策略精选灵活配置混合C','0.08','80.93','3.89','0.0031',4),
    @(11,'009187','天弘聚新三个月定期开放混合C','0.03','23.24','0.99','0.0003',9)
)

# Columns B (fund code), D, E, F, G (numeric-looking text figures) must stay
# text so things like leading zeros ("007202") and trailing zeros ("71.80")
# survive intact instead of being auto-coerced into numbers.
$q1Sheet.Range("B2:B13").NumberFormat = "@"
$q1Sheet.Range("D2:G13").NumberFormat = "@"

$r = 2
foreach ($row in $rows2022Q1) {
    $q1Sheet.Cells.Item($r, 1).Value = $row[0]
    $q1Sheet.Cells.Item($r, 2).Value = $row[1]
    $q1Sheet.Cells.Item($r, 3).Value = $row[2]
    $q1Sheet.Cells.Item($r, 4).Value = $row[3]
    $q1Sheet.Cells.Item($r, 5).Value = $row[4]
    $q1Sheet.Cells.Item($r, 6).Value = $row[5]
    $q1Sheet.Cells.Item($r, 7).Value = $row[6]
    $q1Sheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Header row + column-A index cells use the bold/centered/bordered style
# already used throughout the workbook. Pull it from the untouched
# "2021-Q4" sheet (same B1:H1 / A2:A13 footprint) so every header/index
# cell - old or newly added - ends up on the same shared style index.
$styleSource = $wb.Worksheets.Item(5)
$styleSource.Range("B1").Copy() | Out-Null
$q1Sheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$styleSource.Range("A2").Copy() | Out-Null
$q1Sheet.Range("A2:A13").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet at the end of the workbook holding
# the quarter-over-quarter summary (original rows + new 2022-Q1 row).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$totalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "总计"
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$totalHeaders = @('日期','持有数量(只)','持有市值(亿元)')
for ($c = 0; $c -lt $totalHeaders.Length; $c++) {
    $totalSheet.Cells.Item(1, 2 + $c).Value = $totalHeaders[$c]
}

$rowsTotal = @(
    @(0,'2022-Q1',12,0.97),
    @(1,'2021-Q4',12,1.08),
    @(2,'2021-Q3',8,0.97),
    @(3,'2021-Q2',10,2.17),
    @(4,'2021-Q1',8,5.03),
    @(5,'2020-Q4',5,2.25)
)

$r = 2
foreach ($row in $rowsTotal) {
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Same header/index styling as every other sheet in the workbook.
$styleSource.Range("B1").Copy() | Out-Null
$totalSheet.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$styleSource.Range("A2").Copy() | Out-Null
$totalSheet.Range("A2:A7").PasteSpecial(-4122) | Out-Null

# Re-activate the first sheet so the workbook's active tab stays where it
# was before we added/activated the new trailing sheet.
$wb.Worksheets.Item(1).Activate()

Write-Output "done"
